# Update countries & provincias Spain
# Refresh the COVID-19 country data snapshot: update the "last updated"
# timestamp, swap the rank order of Venezuela/Nigeria and Belice/Liberia
# (their case counts crossed), and refresh the statistics for the rows
# whose underlying source data changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp row (A1) ---------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 9 de Septiembre de 2020 a las 06:46"

# --- India (row 5) ------------------------------------------------------
$ws.Range("B5").Value = 4370128
$ws.Range("C5").Value = 2692
$ws.Range("D5").Value = 3398844
$ws.Range("E5").Value = 897361

# --- Pakistan (row 20) ---------------------------------------------------
$ws.Range("B20").Value = 299659
$ws.Range("C20").Value = 426
$ws.Range("D20").Value = 286506
$ws.Range("E20").Value = 6794
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 6359

# --- Honduras (row 50) ----------------------------------------------------
$ws.Range("B50").Value = 65218
$ws.Range("C50").Value = 404
$ws.Range("D50").Value = 14273
$ws.Range("E50").Value = 48911
$ws.Range("G50").Value = 11
$ws.Range("H50").Value = 2034

# --- Venezuela / Nigeria swap (rows 55 & 56) ------------------------------
# Venezuela overtook Nigeria in total cases, so it now ranks above Nigeria.
$ws.Range("A55").Value = "Venezuela"
$ws.Range("B55").Value = 55563
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 44435
$ws.Range("E55").Value = 10684
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 444

$ws.Range("A56").Value = "Nigeria"
$ws.Range("B56").Value = 55456
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 43334
$ws.Range("E56").Value = 11055
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 1067

# --- Tailandia (row 125) --------------------------------------------------
$ws.Range("B125").Value = 3447
$ws.Range("C125").Value = 1
$ws.Range("D125").Value = 3286
$ws.Range("E125").Value = 103

# --- Belice / Liberia swap (rows 160 & 161) -------------------------------
# Belice overtook Liberia in total cases, so it now ranks above Liberia.
$ws.Range("A160").Value = "Belice"
$ws.Range("B160").Value = 1361
$ws.Range("C160").Value = 54
$ws.Range("D160").Value = 321
$ws.Range("E160").Value = 1024
$ws.Range("F160").Value = 0
$ws.Range("G160").Value = 9
$ws.Range("H160").Value = 16

$ws.Range("A161").Value = "Liberia"
$ws.Range("B161").Value = 1311
$ws.Range("C161").Value = 0
$ws.Range("D161").Value = 1194
$ws.Range("E161").Value = 35
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 0
$ws.Range("H161").Value = 82

# --- Butan (row 187) -------------------------------------------------------
$ws.Range("B187").Value = 234
$ws.Range("C187").Value = 1
$ws.Range("D187").Value = 153
$ws.Range("E187").Value = 81
